$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=4.279621666666666; H=12.838865; I=0.07738101347700441; J=0.07738101347700439; K=3; M=17.16653; N=51.49959; O=0.0560345397128279; P=0.0560345397128279; Q=73.46625372948333; R=661.1962835653499; S=0.004336009472696074; T=0.004336009472696073 }
    3  = @{ E=3; G=4.279621666666666; H=12.838865; I=0.07738101347700441; J=0.07738101347700439; K=3; M=256.4443053333333; N=769.332916; O=0.8370788162388805; P=0.8370788162388805; Q=1097.484605397815; R=9877.361448580337; S=0.0647740071606957; T=0.06477400716069569 }
    4  = @{ E=3; G=4.279621666666666; H=12.838865; I=0.07738101347700441; J=0.07738101347700439; K=3; M=32.74538866666666; N=98.236166; O=0.1068866440482915; P=0.1068866440482915; Q=140.1378748212878; R=1261.24087339159; S=0.00827099684361262; T=0.008270996843612618 }
    5  = @{ E=3; G=24.807588; H=74.422764; I=0.4485528046349828; J=0.4485528046349828; K=3; M=17.16653; N=51.49959; O=0.0560345397128279; P=0.0560345397128279; Q=425.86020362964; R=3832.74183266676; S=0.02513444994461928; T=0.02513444994461928 }
    6  = @{ E=3; G=24.807588; H=74.422764; I=0.4485528046349828; J=0.4485528046349828; K=3; M=256.4443053333333; N=769.332916; O=0.8370788162388805; P=0.8370788162388805; Q=6361.764671655535; R=57255.88204489982; S=0.3754740507244813; T=0.3754740507244813 }
    7  = @{ E=3; G=24.807588; H=74.422764; I=0.4485528046349828; J=0.4485528046349828; K=3; M=32.74538866666666; N=98.236166; O=0.1068866440482915; P=0.1068866440482915; Q=812.3341109425359; R=7311.006998482824; S=0.04794430396588226; T=0.04794430396588226 }
    8  = @{ E=3; G=26.21862666666667; H=78.65588; I=0.4740661818880128; J=0.4740661818880128; K=3; M=17.16653; N=51.49959; O=0.0560345397128279; P=0.0560345397128279; Q=450.0828412321333; R=4050.7455710892; S=0.02656408029551255; T=0.02656408029551255 }
    9  = @{ E=3; G=26.21862666666667; H=78.65588; I=0.4740661818880128; J=0.4740661818880128; K=3; M=256.4443053333333; N=769.332916; O=0.8370788162388805; P=0.8370788162388805; Q=6723.617502327342; R=60512.55752094607; S=0.3968307583537036; T=0.3968307583537036 }
    10 = @{ E=3; G=26.21862666666667; H=78.65588; I=0.4740661818880128; J=0.4740661818880128; K=3; M=32.74538866666666; N=98.236166; O=0.1068866440482915; P=0.1068866440482915; Q=858.539120506231; R=7726.85208455608; S=0.05067134323879665; T=0.05067134323879665 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
